# QB Website updated 12/13
# Update the "Occurrence" values for the first two question rows to append
# the newly-added occurrence numbers (83, 87), which grows a new shared
# string and bumps row 3's height to fit the extra text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newOccurrence = "25, 27, 31, 35, 39, 47, 51, 55, 59, 63, 67, 69, 75, 83, 87"

$ws.Range("E2").Value = $newOccurrence
$ws.Range("E3").Value = $newOccurrence

# Row 3 needs to grow to fit the longer wrapped text.
$ws.Rows.Item(3).RowHeight = 46

# Reflect the updated on-screen selection/scroll position.
$ws.Range("G6").Select()
